$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 938 - this shifts the existing rows
# 938-966 down to 941-969 (Excel also extends the used range / dimension
# automatically, and copies formatting, e.g. the date style on column D,
# from the row above into the freshly inserted rows).
$ws.Rows.Item(938).Resize(3).Insert()

# Common (constant) columns shared by every "Frutilla" row in this block.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$codreg = 4
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

# New weekly report rows for the date 2023-01-13 (Excel serial 44939):
# Especial, Primera and Segunda quality records.
$newRows = @(
    @{ Row = 938; Calidad = "Especial"; Volumen = 400; PMin = 13000; PMax = 14000; PProm = 13500; PKg = 1929 },
    @{ Row = 939; Calidad = "Primera";  Volumen = 500; PMin = 11000; PMax = 12000; PProm = 11500; PKg = 1643 },
    @{ Row = 940; Calidad = "Segunda";  Volumen = 500; PMin = 9000;  PMax = 10000; PProm = 9500;  PKg = 1357 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = 44939
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
